$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the header row (row 1) to reflect the "real data" test run:
#   A1: z_real     -> Z_real
#   B1: z_imag     -> Z_imag
#   C1: frequency  -> angular frequency
# (D1:G1 keep their existing labels: eff_cap+D1:D37, applied voltage, J_ph, J)
#
# Assignment order matches the order new labels were introduced into the
# workbook's shared-string table (angular frequency, then Z_real, then Z_imag).
$ws.Range("C1").Value = "angular frequency"
$ws.Range("A1").Value = "Z_real"
$ws.Range("B1").Value = "Z_imag"

# Move the active cell selection to B1
$ws.Range("B1").Select()
